$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last row's order count, revenue and cost figures.
$ws.Range("B8").Value = 700
$ws.Range("C8").Value = 50000
$ws.Range("D8").Value = 30000

# Re-enter the gross-profit formula for row 8 as its own (non-shared) formula.
$ws.Range("E8").Formula = "=C8-D8"

# Move the active selection to B7, matching the saved view state.
$ws.Range("B7").Select()
